# Add one new data row to each of the 4 worksheets (the sheets are log
# tables of hex/numeric records). Each sheet gets exactly one new row
# appended right after its current last row.

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $timeVal, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $gIsText, $hVal, $iVal) {
    # Column A: date/time serial, formatted the same way as the existing
    # rows above it (numFmt "YYYY-MM-DD HH:MM:SS").
    $ws.Cells.Item($row, 1).Value = $timeVal
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B-E: plain text values (stored as inline/shared strings).
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal

    # Column F: plain number.
    $ws.Cells.Item($row, 6).Value = $fVal

    # Column G: either a plain number, or (when too large to round-trip
    # exactly / kept as text in the source) a text value, restored to the
    # default (General) cell style afterwards so formatting matches the
    # surrounding, un-styled cells.
    if ($gIsText) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $gVal
        $ws.Cells.Item($row, 7).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 7).Value = $gVal
    }

    # Columns H-I: plain numbers.
    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

$gNum = [double]"5.68631262647114e+23"

# --- Sheet 1: ROW50-FE-LIFTER -> append row 60 ---
$ws1 = $wb.Worksheets.Item(1)
Add-LogRow $ws1 60 45754.21704342592 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x5a" "0xe" 400 $gNum $false 346 14

# --- Sheet 2: ROW50-MID-LIFTER -> append row 62 ---
$ws2 = $wb.Worksheets.Item(2)
Add-LogRow $ws2 62 45754.18472222222 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x5e" "0x19" 400 "568631262647113771663628" $true 350 25

# --- Sheet 3: ROW11-FE-LIFTER -> append row 60 ---
$ws3 = $wb.Worksheets.Item(3)
Add-LogRow $ws3 60 45754.25228775463 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x5a" "0x14" 400 $gNum $false 346 20

# --- Sheet 4: ROW11-MID-LIFTER -> append row 60 ---
$ws4 = $wb.Worksheets.Item(4)
Add-LogRow $ws4 60 45754.38132809028 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x62" "0x19" 400 $gNum $false 354 25

Write-Host "Added new rows to all 4 sheets"
